$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M3").Value = 1.22
$ws.Range("N3").Value = 4
$ws.Range("Q3").Value = 3.15
$ws.Range("R3").Value = 1.37
$ws.Range("U3").Value = 8.199999999999999
$ws.Range("W3").Value = 11
$ws.Range("X3").Value = 1.05
$ws.Range("Y3").Value = 1.95
$ws.Range("Z3").Value = 1.85
$ws.Range("G4").Value = 3.8
$ws.Range("I4").Value = 2.2
$ws.Range("L4").Value = 3.1
$ws.Range("AC4").Value = 7.5
$ws.Range("AP4").Value = 21
$ws.Range("G6").Value = 2.5
$ws.Range("I6").Value = 3.4
$ws.Range("J6").Value = 3.4
$ws.Range("L6").Value = 4.33
$ws.Range("AC6").Value = 5.5
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 11
$ws.Range("AF6").Value = 23
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 41
$ws.Range("AI6").Value = 5.5
$ws.Range("AJ6").Value = 6
$ws.Range("AL6").Value = 101
$ws.Range("AM6").Value = 7
$ws.Range("AN6").Value = 15
$ws.Range("AP6").Value = 41
$ws.Range("G8").Value = 1.85
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("AS8").Value = 1250
$ws.Range("G9").Value = 3.25
$ws.Range("I9").Value = 2.4
$ws.Range("AG9").Value = 34
$ws.Range("AM9").Value = 6
$ws.Range("AO9").Value = 10
$ws.Range("H10").Value = 2.57
$ws.Range("I10").Value = 3.1
$ws.Range("J10").Value = 3.45
$ws.Range("L10").Value = 3.95
$ws.Range("N10").Value = 4.45
$ws.Range("O10").Value = 1.65
$ws.Range("P10").Value = 2.1
$ws.Range("S10").Value = 2.95
$ws.Range("T10").Value = 1.35
$ws.Range("W10").Value = 5.4
$ws.Range("X10").Value = 1.12
$ws.Range("Y10").Value = 1.65
$ws.Range("Z10").Value = 2.1
$ws.Range("AA10").Value = 2.25
$ws.Range("AC10").Value = 5.8
$ws.Range("AD10").Value = 11.75
$ws.Range("AE10").Value = 10.75
$ws.Range("AG10").Value = 30
$ws.Range("AH10").Value = 55
$ws.Range("AI10").Value = 4.45
$ws.Range("AM10").Value = 6.1
$ws.Range("AN10").Value = 14
$ws.Range("AO10").Value = 12.5
$ws.Range("AQ10").Value = 40
$ws.Range("AR10").Value = 65
$ws.Range("G11").Value = 2.62
$ws.Range("H11").Value = 2.6
$ws.Range("I11").Value = 3.2
$ws.Range("J11").Value = 3.5
$ws.Range("L11").Value = 4.05
$ws.Range("M11").Value = 1.18
$ws.Range("X11").Value = 1.09
$ws.Range("Y11").Value = 1.72
$ws.Range("Z11").Value = 2
$ws.Range("AA11").Value = 2.42
$ws.Range("AC11").Value = 5.2
$ws.Range("AD11").Value = 10.75
$ws.Range("AE11").Value = 11.75
$ws.Range("AF11").Value = 32
$ws.Range("AG11").Value = 35
$ws.Range("AJ11").Value = 5.6
$ws.Range("AK11").Value = 23
$ws.Range("AL11").Value = 200
$ws.Range("AO11").Value = 13
$ws.Range("AQ11").Value = 45
$ws.Range("AR11").Value = 75
$ws.Range("G15").Value = 4.25
$ws.Range("H15").Value = 3.8
$ws.Range("I15").Value = 1.7
$ws.Range("J15").Value = 4.5
$ws.Range("L15").Value = 2.22
$ws.Range("T15").Value = 1.93
$ws.Range("AE15").Value = 14
$ws.Range("AH15").Value = 40
$ws.Range("AI15").Value = 12
$ws.Range("AM15").Value = 7.6
$ws.Range("AP15").Value = 13
$ws.Range("M16").Value = 1.05
$ws.Range("N16").Value = 8.5
$ws.Range("S16").Value = 1.93
$ws.Range("T16").Value = 1.88
$ws.Range("H17").Value = 3.5
$ws.Range("I17").Value = 1.95
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 9
$ws.Range("O17").Value = 1.22
$ws.Range("P17").Value = 4
$ws.Range("S17").Value = 1.75
$ws.Range("T17").Value = 2.05
$ws.Range("W17").Value = 2.75
$ws.Range("X17").Value = 1.4
$ws.Range("AA17").Value = 1.67
$ws.Range("AB17").Value = 2.1
$ws.Range("AG17").Value = 26
$ws.Range("AI17").Value = 12
$ws.Range("AR17").Value = 23
$ws.Range("G18").Value = 2.63
$ws.Range("I18").Value = 2.45
$ws.Range("M18").Value = 1.03
$ws.Range("N18").Value = 10
$ws.Range("S18").Value = 1.6
$ws.Range("T18").Value = 2.3
$ws.Range("U18").Value = 1.98
$ws.Range("V18").Value = 1.83
$ws.Range("W18").Value = 2.38
$ws.Range("X18").Value = 1.53
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 3.6
$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 9.5
$ws.Range("O19").Value = 1.3
$ws.Range("P19").Value = 3.4
$ws.Range("S19").Value = 2.03
$ws.Range("T19").Value = 1.78
$ws.Range("W19").Value = 3.5
$ws.Range("X19").Value = 1.29
$ws.Range("AF19").Value = 17
$ws.Range("AI19").Value = 9.5
$ws.Range("AM19").Value = 10
$ws.Range("AN19").Value = 19
$ws.Range("G21").Value = 2.05
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 3.3
$ws.Range("J21").Value = 2.75
$ws.Range("K21").Value = 2.2
$ws.Range("L21").Value = 4
$ws.Range("O21").Value = 1.29
$ws.Range("P21").Value = 3.5
$ws.Range("S21").Value = 1.95
$ws.Range("T21").Value = 1.85
$ws.Range("W21").Value = 3.4
$ws.Range("X21").Value = 1.3
$ws.Range("AA21").Value = 1.73
$ws.Range("AB21").Value = 2
$ws.Range("AC21").Value = 8
$ws.Range("AD21").Value = 10
$ws.Range("AF21").Value = 19
$ws.Range("AG21").Value = 17
$ws.Range("AH21").Value = 26
$ws.Range("AS21").Value = 201
$ws.Range("N22").Value = 13
$ws.Range("G23").Value = 1.83
$ws.Range("H23").Value = 3.9
$ws.Range("I23").Value = 3.6
$ws.Range("J23").Value = 2.3
$ws.Range("K23").Value = 2.4
$ws.Range("L23").Value = 3.9
$ws.Range("M23").Value = 1.03
$ws.Range("N23").Value = 9.25
$ws.Range("O23").Value = 1.16
$ws.Range("P23").Value = 4.65
$ws.Range("S23").Value = 1.5
$ws.Range("T23").Value = 2.42
$ws.Range("W23").Value = 2.15
$ws.Range("X23").Value = 1.62
$ws.Range("Y23").Value = 1.27
$ws.Range("Z23").Value = 3.4
$ws.Range("AA23").Value = 1.47
$ws.Range("AB23").Value = 2.5
$ws.Range("AC23").Value = 11.25
$ws.Range("AD23").Value = 11.5
$ws.Range("AE23").Value = 8.5
$ws.Range("AF23").Value = 17
$ws.Range("AG23").Value = 12.5
$ws.Range("AH23").Value = 17.5
$ws.Range("AI23").Value = 9.25
$ws.Range("AJ23").Value = 8
$ws.Range("AL23").Value = 37
$ws.Range("AN23").Value = 24
$ws.Range("AO23").Value = 12.5
$ws.Range("AP23").Value = 55
$ws.Range("AQ23").Value = 27
$ws.Range("AR23").Value = 27
$ws.Range("AS23").Value = 200
$ws.Range("G24").Value = 1.52
$ws.Range("H24").Value = 4.35
$ws.Range("I24").Value = 5.2
$ws.Range("J24").Value = 1.98
$ws.Range("K24").Value = 2.47
$ws.Range("M24").Value = 1.03
$ws.Range("N24").Value = 9.25
$ws.Range("O24").Value = 1.16
$ws.Range("P24").Value = 4.55
$ws.Range("S24").Value = 1.5
$ws.Range("T24").Value = 2.42
$ws.Range("W24").Value = 2.2
$ws.Range("X24").Value = 1.6
$ws.Range("Y24").Value = 1.27
$ws.Range("Z24").Value = 3.4
$ws.Range("AA24").Value = 1.6
$ws.Range("AB24").Value = 2.18
$ws.Range("AC24").Value = 9.5
$ws.Range("AD24").Value = 8.75
$ws.Range("AE24").Value = 8.25
$ws.Range("AF24").Value = 11.75
$ws.Range("AG24").Value = 11
$ws.Range("AH24").Value = 19.5
$ws.Range("AI24").Value = 9.25
$ws.Range("AJ24").Value = 8.75
$ws.Range("AK24").Value = 14.5
$ws.Range("AL24").Value = 50
$ws.Range("AM24").Value = 19
$ws.Range("AN24").Value = 35
$ws.Range("AO24").Value = 16.5
$ws.Range("AP24").Value = 90
$ws.Range("AQ24").Value = 45
$ws.Range("AR24").Value = 40
$ws.Range("AS24").Value = 300
